$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.283.26'
$ws.Range('E2').Value = '  +3.66%  '
$ws.Range('D3').Value = '3.206.10'
$ws.Range('E3').Value = '  +2.13%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '539.86'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.29'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +4.31%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.529'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +3.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.36'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.58%  '
$ws.Range('E10').Value = '  +4.21%  '
$ws.Range('E11').Value = '  +1.23%  '
$ws.Range('D12').Value = '3.753.58'
$ws.Range('E12').Value = '  +2.05%  '
$ws.Range('E13').Value = '  -1.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000176'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +3.75%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.16'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.43%  '
$ws.Range('D16').Value = '60.286.61'
$ws.Range('E16').Value = '  +3.48%  '
$ws.Range('D17').Value = '3.171.54'
$ws.Range('E17').Value = '  +1.17%  '
$ws.Range('E18').Value = '  -0.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.14'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.16%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.38'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +2.32%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '383.75'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +2.35%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.531'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.90%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.32'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.18%  '
$ws.Range('E25').Value = '  +2.38%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.81'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +9.22%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('D28').Value = '0.0₃0910'
$ws.Range('E28').Value = '  +2.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.92'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.34%  '
$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.46'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +5.74%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.43'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +2.95%  '
$ws.Range('B32').Value = 'RenderToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.20'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.57%  '
$ws.Range('E33').Value = '  +4.77%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.60'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +5.26%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '156.77'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -3.31%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.37'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.20%  '
$ws.Range('B37').Value = 'EnergySwap'
$ws.Range('C37').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '25.88'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.38%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').Value = '2.777.01'
$ws.Range('E38').Value = '  +4.86%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0711'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +4.59%  '
$ws.Range('E40').Value = '  +0.25%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.24'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.33%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.80'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +2.87%  '
$ws.Range('E43').Value = '  +3.79%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0286'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +4.15%  '
$ws.Range('D45').Value = '3.245.84'
$ws.Range('E45').Value = '  +2.05%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.01'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +3.22%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.102'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.53%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.16'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.18%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.806'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +7.62%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.60'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.44%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.999'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.04%  '
